$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Booking Facility Name" / "Booking Facility Number" rows (no longer mapped)
$ws.Rows("17:18").Delete()

# Update the Consent Decision Code comment cell (now shifted up to row 21) to
# summarize the two remaining codes instead of listing each one in its own row
$ws.Range("B21").Value = "Codes: Consent Granted; Consent Denied"

# Remove the now-obsolete individual code rows ("Consent Denied", "Inmate Never
# Seen", "Consent Not Obtained")
$ws.Rows("22:24").Delete()

# Keep the previously-selected cell reference on the bottom-right pane
$ws.Range("B29").Select()
